$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 1010
$ws.Range("B5").Value = "STU"
$ws.Range("C5").Value = "VWX"

$ws.Range("D6").Select()
